$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.119.43"
$ws.Range("E2").Value = "  -1.07%  "
$ws.Range("D3").Value = "1.846.00"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'236.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.4768"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").Value = "'0.2807"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.47%  "
$ws.Range("D9").Value = "'0.06469"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("D10").Value = "1.858.40"
$ws.Range("E10").Value = "  -1.75%  "
$ws.Range("D11").Value = "'0.07314"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "'16.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.46%  "
$ws.Range("D13").Value = "'5.110"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "'87.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("D15").Value = "'0.6445"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").Value = "30.070.27"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.02%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").Value = "'0.000007619"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.65%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.099.09"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'222.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +16.92%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'5.288"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("D24").Value = "'6.079"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").Value = "'9.213"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.00%  "
$ws.Range("D26").Value = "'163.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("D27").Value = "'18.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.03%  "
$ws.Range("D28").Value = "'1.914"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "'0.09201"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("D31").Value = "'4.235"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("D32").Value = "'3.957"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "'0.05007"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.03%  "
$ws.Range("D34").Value = "'0.7381"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.08%  "
$ws.Range("D36").Value = "'2.688"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'0.01818"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'2.603"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.9043"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.051"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.933"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'106.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "'0.4247"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").Value = "'0.9991"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "'7.415"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("D46").Value = "'0.1312"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.76%  "
$ws.Range("D47").Value = "'1.560"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.85%  "
$ws.Range("D48").Value = "'63.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.87%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.768"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").Value = "'0.05666"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.63%  "
